# Auto-generated Excel COM-interop script
# Applies the 2024-05-07 daily violent-crime data update
# to the "violent-crime-full-year.xlsx" workbook.

$wb = $excel.ActiveWorkbook

$changes = @{
    'Citywide Totals' = @{ 'K2' = '2528'; 'K3' = '2449'; 'G4' = '1483'; 'J4' = '1816'; 'K4' = '512'; 'K5' = '162'; 'K6' = '3048'; 'G7' = '24709'; 'J7' = '29283'; 'K7' = '8699' }
    'Logan Square' = @{ 'K6' = '67'; 'K7' = '128' }
    'Austin' = @{ 'K2' = '169'; 'K3' = '174'; 'K6' = '193'; 'K7' = '580' }
    'South Chicago' = @{ 'K6' = '44'; 'K7' = '192' }
    'Garfield Park' = @{ 'K2' = '97'; 'K3' = '124'; 'K6' = '94'; 'K7' = '343' }
    'Grand Crossing' = @{ 'K3' = '97'; 'K6' = '85'; 'K7' = '279' }
    'New City' = @{ 'K6' = '83'; 'K7' = '205' }
    'Woodlawn' = @{ 'K4' = '7'; 'K6' = '43'; 'K7' = '155' }
    'Fuller Park' = @{ 'K6' = '9'; 'K7' = '25' }
    'By Neighborhood' = @{ 'K2' = '62'; 'K4' = '33'; 'K6' = '71'; 'K7' = '257'; 'K8' = '580'; 'K14' = '51'; 'K15' = '87'; 'K16' = '26'; 'K18' = '58'; 'K19' = '255'; 'K20' = '198'; 'K21' = '23'; 'K23' = '78'; 'K29' = '441'; 'K30' = '25'; 'K31' = '99'; 'K33' = '343'; 'K37' = '279'; 'K42' = '305'; 'K44' = '84'; 'K47' = '48'; 'K48' = '109'; 'K51' = '95'; 'K52' = '237'; 'K53' = '128'; 'K54' = '162'; 'K60' = '57'; 'G63' = '284'; 'K63' = '34'; 'K64' = '56'; 'K65' = '205'; 'K66' = '30'; 'K67' = '341'; 'K70' = '16'; 'K72' = '41'; 'K74' = '11'; 'K76' = '126'; 'K78' = '119'; 'K83' = '192'; 'K84' = '63'; 'K85' = '418'; 'K86' = '55'; 'K89' = '114'; 'J91' = '338'; 'K91' = '83'; 'K99' = '155'; 'G101' = '24709'; 'J101' = '29283'; 'K101' = '8699' }
    'Gage Park' = @{ 'K6' = '39'; 'K7' = '99' }
    'North Lawndale' = @{ 'K3' = '106'; 'K7' = '341' }
    'South Deering' = @{ 'K3' = '20'; 'K7' = '63' }
    'Loop' = @{ 'K4' = '9'; 'K7' = '162' }
    'Englewood' = @{ 'K2' = '119'; 'K3' = '147'; 'K4' = '25'; 'K7' = '441' }
    'Lake View' = @{ 'K3' = '20'; 'K7' = '109' }
    'Chatham' = @{ 'K2' = '80'; 'K3' = '68'; 'K7' = '255' }
    'Irving Park' = @{ 'K3' = '23'; 'K7' = '84' }
    'River North' = @{ 'K6' = '77'; 'K7' = '126' }
    'Bridgeport' = @{ 'K2' = '21'; 'K3' = '10'; 'K7' = '51' }
    'Ashburn' = @{ 'K2' = '25'; 'K7' = '71' }
    'Humboldt Park' = @{ 'K3' = '95'; 'K7' = '305' }
    'Rogers Park' = @{ 'K3' = '28'; 'K6' = '45'; 'K7' = '119' }
    'Douglas' = @{ 'K3' = '23'; 'K7' = '78' }
    'Washington Park' = @{ 'K3' = '35'; 'J4' = '17'; 'J7' = '338'; 'K7' = '83' }
    'Chinatown' = @{ 'K6' = '13'; 'K7' = '23' }
    'Near South Side' = @{ 'K6' = '19'; 'K7' = '56' }
    'Chicago Lawn' = @{ 'K3' = '56'; 'K6' = '68'; 'K7' = '198' }
    'Calumet Heights' = @{ 'K3' = '19'; 'K5' = '2'; 'K7' = '58' }
    'Auburn Gresham' = @{ 'K3' = '82'; 'K6' = '69'; 'K7' = '257' }
    'Kenwood' = @{ 'K6' = '15'; 'K7' = '48' }
    'Brighton Park' = @{ 'K2' = '30'; 'K7' = '87' }
    'North Center' = @{ 'K3' = '5'; 'K7' = '30' }
    'Albany Park' = @{ 'K3' = '17'; 'K6' = '23'; 'K7' = '62' }
    'O''Hare' = @{ 'K2' = '6'; 'K7' = '16' }
    'Uptown' = @{ 'K2' = '27'; 'K7' = '114' }
    'Streeterville' = @{ 'K3' = '12'; 'K7' = '55' }
    'Little Italy, UIC' = @{ 'K2' = '25'; 'K7' = '95' }
    'Morgan Park' = @{ 'K2' = '16'; 'K7' = '57' }
    'South Shore' = @{ 'K2' = '155'; 'K3' = '142'; 'K7' = '418' }
    'Old Town' = @{ 'K6' = '24'; 'K7' = '41' }
    'Little Village' = @{ 'K6' = '99'; 'K7' = '237' }
    'Archer Heights' = @{ 'K3' = '10'; 'K7' = '33' }
    'Bucktown' = @{ 'K2' = '6'; 'K7' = '26' }
    'Printers Row' = @{ 'K6' = '8'; 'K7' = '11' }
}

foreach ($sheetName in $changes.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $cellChanges = $changes[$sheetName]
    foreach ($cellRef in $cellChanges.Keys) {
        $ws.Range($cellRef).Value = [double]$cellChanges[$cellRef]
    }
}

Write-Output "Applied $($changes.Keys.Count) sheet updates."
